$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-28 Sunday" "2025-09-29 Monday"

Replace-Text "533÷3=" "305÷4="
Replace-Text "233÷2=" "859÷8="
Replace-Text "647÷9=" "315÷9="
Replace-Text "168÷7=" "249÷4="
Replace-Text "603÷7=" "963÷7="

Replace-Text "888÷9=" "764÷8="
Replace-Text "695÷7=" "566÷7="
Replace-Text "192÷4=" "809÷2="
Replace-Text "686÷7=" "365÷7="
Replace-Text "799÷9=" "234÷6="

Replace-Text "581÷6=" "242÷6="
Replace-Text "780÷5=" "651÷7="
Replace-Text "935÷6=" "240÷4="
Replace-Text "591÷7=" "103÷7="
Replace-Text "873÷3=" "463÷6="

Replace-Text "877÷5=" "695÷5="
Replace-Text "814÷7=" "837÷3="
Replace-Text "923÷3=" "581÷2="
Replace-Text "728÷9=" "340÷7="
Replace-Text "833÷4=" "142÷6="

Replace-Text "937÷4=" "869÷4="
Replace-Text "358÷9=" "693÷6="
Replace-Text "256÷5=" "678÷7="
Replace-Text "164÷7=" "210÷5="
Replace-Text "316÷5=" "899÷5="
